$wb = $excel.ActiveWorkbook

$proximity = $wb.Worksheets.Item("Proximity")
$camera = $wb.Worksheets.Item("Camera")

$newRows = @(
    @{ Sheet = $proximity; Row = 23; Date = "2026-02-01"; Timestamp = "18:16:20"; Hour = "18:00"; Location = "Living Room Main Door"; Value = "ENTER"; Status = "User ENTERED Living Room Main Door" },
    @{ Sheet = $proximity; Row = 24; Date = "2026-02-01"; Timestamp = "18:16:20"; Hour = "18:00"; Location = "Living Room Main Door"; Value = "EXIT";  Status = "User EXITED Living Room Main Door" },
    @{ Sheet = $proximity; Row = 25; Date = "2026-02-01"; Timestamp = "18:16:28"; Hour = "18:00"; Location = "Living Room Main Door"; Value = "ENTER"; Status = "User ENTERED Living Room Main Door" },
    @{ Sheet = $camera;    Row = 23; Date = "2026-02-01"; Timestamp = "18:16:20"; Hour = "18:00"; Location = "Living Room Main Door"; Value = "Image Captured"; Status = "Active" },
    @{ Sheet = $camera;    Row = 24; Date = "2026-02-01"; Timestamp = "18:16:29"; Hour = "18:00"; Location = "Living Room Main Door"; Value = "Image Captured"; Status = "Active" }
)

foreach ($r in $newRows) {
    $ws = $r.Sheet
    $row = $r.Row

    # Column A holds a plain "YYYY-MM-DD" text value. Assigning that string
    # straight to .Value would make Excel auto-convert it into a real date
    # serial (and stamp a date NumberFormat on the cell). Force the cell to
    # Text first so the literal string is preserved, then restore the
    # default "Normal" style so no stray formatting is left on the cell.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 1).Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $r.Timestamp
    $ws.Cells.Item($row, 3).Value = $r.Hour
    $ws.Cells.Item($row, 4).Value = $r.Location
    $ws.Cells.Item($row, 5).Value = $r.Value
    $ws.Cells.Item($row, 6).Value = $r.Status
}
